# Reproduce the "added the transfer diagram for the carrier model" commit:
#  - The slide that already carries the Rectangle/background box, the SIR
#    diagram shapes and the "Disease states" caption (originally slide #2)
#    is moved in front of the bare SIR-diagram slide (originally slide #1),
#    becoming the deck's first slide.
#  - On that (now first) slide, the background "Rectangle 2" shape is
#    nudged slightly (a few EMU) to its final resting position.

$p = $ppt.ActivePresentation

# Identify the slide that holds the "Rectangle 2" shape (the fuller
# "transfer diagram" slide) regardless of its current index.
$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        if ($slide.Shapes.Item($j).Name -eq "Rectangle 2") {
            $targetIndex = $i
            break
        }
    }
    if ($targetIndex -ne -1) { break }
}

# Move that slide to the front of the deck (position 1) if it isn't there.
if ($targetIndex -ne 1) {
    $p.Slides.Item($targetIndex).MoveTo(1)
}

# Nudge the "Rectangle 2" shape on the (now first) slide to its new offset.
$s = $p.Slides.Item(1)
for ($j = 1; $j -le $s.Shapes.Count; $j++) {
    $sh = $s.Shapes.Item($j)
    if ($sh.Name -eq "Rectangle 2") {
        $sh.Left = 300.5066141732283
        $sh.Top = 22.35795275590551
    }
}
